$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the D1/E1 header labels (bedrooms_2 <-> kitchens_2)
$ws.Range("D1").Value = "kitchens_2"
$ws.Range("E1").Value = "bedrooms_2"

# Update the block-order indicator values to match the relabeled columns
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 1

$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0

$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 0
